$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 180 (shifting existing rows 180+ down by one)
$ws.Rows.Item(180).Insert()

# Populate the newly inserted row 180 with the new weekly record
$ws.Cells.Item(180, 1).Value = 3
$ws.Cells.Item(180, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(180, 3).Value = "Coquimbo"
$ws.Cells.Item(180, 4).Value = 44508
$ws.Cells.Item(180, 5).Value = 5
$ws.Cells.Item(180, 6).Value = 100112009
$ws.Cells.Item(180, 7).Value = "Acelga"
$ws.Cells.Item(180, 8).Value = "Sin especificar"
$ws.Cells.Item(180, 9).Value = "Primera"
$ws.Cells.Item(180, 10).Value = 290
$ws.Cells.Item(180, 11).Value = 2000
$ws.Cells.Item(180, 12).Value = 2200
$ws.Cells.Item(180, 13).Value = 2103
$ws.Cells.Item(180, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(180, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(180, 16).Value = 350
$ws.Cells.Item(180, 17).Value = 6
$ws.Cells.Item(180, 18).Value = "Hortaliza"
